$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

$ws.Cells.Item(4, 2).Value = "yes"
$ws.Cells.Item(5, 3).Value = "No "
$ws.Cells.Item(29, 2).Value = "_less_than_20_000"
$ws.Cells.Item(30, 2).Value = "_20_000_to_34_999"
$ws.Cells.Item(48, 2).Value = "prefer_not_to_say"
$ws.Cells.Item(57, 2).Value = "under_6"
$ws.Cells.Item(71, 2).Value = "four_or_more"
$ws.Cells.Item(75, 2).Value = "four_or_more"
$ws.Cells.Item(76, 2).Value = "strongly_disagree"
$ws.Cells.Item(78, 3).Value = "Neutral "
$ws.Cells.Item(80, 3).Value = "strongly agree "
$ws.Cells.Item(83, 2).Value = "i_often_feel_unsafe_when_traveling"
$ws.Cells.Item(84, 2).Value = "usually_takes_too_long"
$ws.Cells.Item(85, 2).Value = "sometimes_takes_too_long"
$ws.Cells.Item(86, 2).Value = "often_go_quickly"
$ws.Cells.Item(87, 2).Value = "always_hard"
$ws.Cells.Item(88, 2).Value = "sometimes_hard"
$ws.Cells.Item(89, 2).Value = "usually_easy"
$ws.Cells.Item(90, 2).Value = "usually_pleasant"
$ws.Cells.Item(91, 2).Value = "sometimes_pleasant"
$ws.Cells.Item(92, 2).Value = "often_unpleasant"
$ws.Cells.Item(93, 2).Value = "usually_afford"
$ws.Cells.Item(94, 2).Value = "sometimes_afford"
$ws.Cells.Item(95, 2).Value = "never_afford"
$ws.Cells.Item(96, 2).Value = "average"
$ws.Cells.Item(97, 2).Value = "smaller"
$ws.Cells.Item(98, 2).Value = "larger"
